# Lab02 presentation fix-up:
#   1. Slide 2 ("Tecnologias Utilizadas") and Slide 3 ("Questões de Pesquisa")
#      had their content swapped in the original file; restore the intended
#      order by moving Slide 2's shapes onto Slide 3 and vice versa.
#   2. Add a new "ck" bullet to the Python technologies list (now living on
#      Slide 3) right before "API GraphQL do GitHub".
#   3. Slide 24's conclusion should read "Não existe relação..." instead of
#      "Existe relação...".

$p = $ppt.ActivePresentation

# --- Part 1: swap the contents of slide 2 and slide 3 -----------------------
$s2 = $p.Slides.Item(2)
$s3 = $p.Slides.Item(3)

# Move all 5 shapes currently on slide 2 (picture, title, body, 2 logo pics)
# to the end of slide 3's shape collection, preserving their relative order.
for ($i = 1; $i -le 5; $i++) {
    $s2.Shapes.Item(1).Cut()
    $s3.Shapes.Paste()
}

# Move the 2 shapes that originally belonged to slide 3 (now sitting at the
# front of slide 3's shape collection) over to slide 2.
for ($i = 1; $i -le 2; $i++) {
    $s3.Shapes.Item(1).Cut()
    $s2.Shapes.Paste()
}

# --- Part 2: insert the new "ck" bullet in the Python tech list ------------
# The Python bullet list now lives on slide 3 as its 3rd shape
# (picture, title, body, pic, pic).
$pySp = $s3.Shapes.Item(3)
$tr = $pySp.TextFrame.TextRange
# "API GraphQL do GitHub" is paragraph 3; insert "ck" as a new paragraph
# right before it, matching its (non-indented) paragraph style.
$apiPara = $tr.Paragraphs(3, 1)
$null = $apiPara.InsertBefore("ck" + [char]13)

# --- Part 3: fix the conclusion wording on slide 24 -------------------------
$s24 = $p.Slides.Item(24)
$shConcl = $s24.Shapes.Item(1)
$trConcl = $shConcl.TextFrame.TextRange
$conclPara = $trConcl.Paragraphs(3, 1)
# Reset first so the engine doesn't try to diff/keep a common run with the
# old text (which would split the replacement across two runs).
$conclPara.Text = "__RESET__"
$conclPara2 = $trConcl.Paragraphs(3, 1)
$conclPara2.Text = "Não existe relação entre a atividade e a qualidade dos repositórios analisados"
